$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# Add two new rows of mobile login test data to the loginTestData sheet
$ws1.Range("A8").Value = "Nho4DGqEoxtXBi1"
$ws1.Range("B8").Value = "RdPS8xmcPF*."
$ws1.Range("C8").Value = "success"

$ws1.Range("A9").Value = "eagermanipulation"
$ws1.Range("B9").Value = "abi3u1nkXd*."
$ws1.Range("C9").Value = "success"

# Move the active selection: googleForgetUserNamePassword sheet is no longer
# the active tab; loginTestData becomes active with C7 selected.
$ws3.Activate()
$ws3.Range("A28").Select()

$ws1.Activate()
$ws1.Range("C7").Select()
